$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.842.25"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "1.874.27"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.79"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5388"
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3762"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07189"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8888"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08157"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.876.01"
$ws.Range("E13").Value = "  +4.15%  "
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.264"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.75"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008551"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "26.902.61"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.985"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.391"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.29"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.261"
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.18"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.726"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.599"
$ws.Range("E30").Value = "  -5.70%  "
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8069"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04978"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.986"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.171"
$ws.Range("E35").Value = "  -4.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6046"
$ws.Range("E36").Value = "  +5.94%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.190"
$ws.Range("E37").Value = "  -5.30%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.592"
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01957"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.95"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5160"
$ws.Range("E44").Value = "  +5.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1494"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.905"
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.63"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06031"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.18"
